# adding averages and more checks
#
# - Restyle the report title and the column-header row on every sheet:
#   bold white text (keeps the existing dark-blue header fill / borders).
# - Training Dashboard sheet: bump the "PERIOD TO EXPIRE" / "LAST UPDATE"
#   values for row 3 (the recalculated check).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # Find the extent of the header row (row 2) so every header cell -
    # including any blank trailing ones - gets the new look.
    $lastCol = $ws.UsedRange.Columns.Count

    $headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $headerRange.Font.Bold = $true
    $headerRange.Font.Color = 16777215   # white

    # Report title (row 1) - same bold/white treatment, normal-sized text.
    $titleCell = $ws.Range("A1")
    $titleCell.Font.Bold = $true
    $titleCell.Font.Size = 11
    $titleCell.Font.Color = 16777215     # white
}

# Training Dashboard: refresh the expiry countdown / last-update check.
$dash = $wb.Worksheets.Item("Training Dashboard")
$dash.Range("H3").Value = 170
$dash.Range("I3").Value = "'16-Sep-2025"
